# Applies the "Add serialization via serialize API #208" restructuring:
#   - drop the empty "Sheet1"
#   - keep "Astruct" as-is (just becomes the active tab)
#   - replace the lone "Bstruct" sheet with four sheets, in this order:
#       Astruct, AstructBstruct2Use, AstructBstructUse, Bstruct, Dstruct
#     where AstructBstruct2Use / AstructBstructUse / Dstruct are brand new,
#     narrow "use" sheets, and Bstruct is re-created with its original data.

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# column width in Excel COM is stored padded by ~5/6 of a character vs. the
# raw OOXML <col width="..."> value - back that out so the saved width
# matches the target exactly.
$padding = 0.8333333333333334

# --- 1. drop Sheet1, purge the old Bstruct sheet (so its shared strings
#        "Floatfield2"/"B2" free up and get re-issued fresh indices later,
#        matching the target shared-string table order) ------------------
[void]$wb.Worksheets.Item("Sheet1").Delete()
[void]$wb.Worksheets.Item("Bstruct").Delete()

$astruct = $wb.Worksheets.Item("Astruct")

# --- 2. AstructBstruct2Use (right after Astruct) -------------------------
$u1 = $wb.Worksheets.Add($null, $astruct)
$u1.Name = "AstructBstruct2Use"
$u1.Range("A1").Value = "Name"
$u1.Range("B1").Value = "Bstrcut2"
$u1.Columns.Item(1).ColumnWidth = 6 - $padding
$u1.Columns.Item(2).ColumnWidth = 10 - $padding
[void]$u1.Range("A1:B1").AutoFilter()

# --- 3. AstructBstructUse (right after AstructBstruct2Use) ---------------
$u2 = $wb.Worksheets.Add($null, $u1)
$u2.Name = "AstructBstructUse"
$u2.Range("A1").Value = "Name"
$u2.Range("B1").Value = "Bstruct2"
$u2.Columns.Item(1).ColumnWidth = 6 - $padding
$u2.Columns.Item(2).ColumnWidth = 10 - $padding
[void]$u2.Range("A1:B1").AutoFilter()

# --- 4. Bstruct, re-created with its original data (right after
#        AstructBstructUse) ----------------------------------------------
$b = $wb.Worksheets.Add($null, $u2)
$b.Name = "Bstruct"
$b.Range("A1").Value = "Name"
$b.Range("B1").Value = "Floatfield"
$b.Range("C1").Value = "Floatfield2"
$b.Range("D1").Value = "Intfield"
$b.Range("A2").Value = "B1"
$b.Range("B2").Value = "0.000000"
$b.Range("C2").Value = "0.000000"
$b.Range("D2").Value = "0"
$b.Range("A3").Value = "B2"
$b.Range("B3").Value = "0.000000"
$b.Range("C3").Value = "0.000000"
$b.Range("D3").Value = "0"
$b.Columns.Item(1).ColumnWidth = 6 - $padding
$b.Columns.Item(2).ColumnWidth = 12 - $padding
$b.Columns.Item(3).ColumnWidth = 13 - $padding
$b.Columns.Item(4).ColumnWidth = 10 - $padding
[void]$b.Range("A1:D1").AutoFilter()

# --- 5. Dstruct (right after Bstruct, last tab) ---------------------------
$d = $wb.Worksheets.Add($null, $b)
$d.Name = "Dstruct"
$d.Range("A1").Value = "Name"
$d.Columns.Item(1).ColumnWidth = 6 - $padding
[void]$d.Range("A1:A1").AutoFilter()

# --- 6. Astruct is the selected/active tab (it inherited that role from
#        the now-deleted Sheet1) ------------------------------------------
$astruct.Activate()
